$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Books")

# Update the Year/Price values to remove ambiguity of the "value" validation
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 9
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 15.65
$ws.Range("C4").Value = -125

# Update the active selection to reflect where the user ended up
$ws.Range("F9").Select()
